$wb = $excel.ActiveWorkbook

$wsSummary   = $wb.Worksheets.Item("Summary")
$wsSchedule  = $wb.Worksheets.Item("Repayment schedule")
$wsTrans     = $wb.Worksheets.Item("Transactions")

# --- Summary sheet: waived penalty recalculation ---
# Row 5 totals drop from 26.64 to 17.76 (one 8.88 installment penalty waived)
$wsSummary.Range("A5").Value = 17.76
$wsSummary.Range("E5").Value = 17.76
$wsSummary.Range("F5").Value = 17.76

# Row 6 (all zero placeholder row) is no longer needed
$wsSummary.Rows("6:6").Delete()

# Update the remembered selection on this sheet
[void]$wsSummary.Range("F5").Select()

# --- Repayment schedule sheet: drop the now-unused duplicate "Over Due" column (O) ---
$wsSchedule.Range("P2").Clear()
$wsSchedule.Range("O3").Clear()
$wsSchedule.Range("O4").Clear()

# Installment 3 penalty (J5) is waived, rolling K5/P5 down to 887.72
$wsSchedule.Range("J5").Value = 0
$wsSchedule.Range("K5").Value = 887.72
$wsSchedule.Range("O5").Clear()
$wsSchedule.Range("P5").Value = 887.72

$wsSchedule.Range("O6").Clear()
$wsSchedule.Range("O7").Clear()
$wsSchedule.Range("O8").Clear()

# Update the remembered selection on this sheet
[void]$wsSchedule.Range("J3:J4").Select()

# --- Switch the active tab from NewLoanInput to Transactions ---
$wsTrans.Activate()
